$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "CasesTab" query (row 2, column B) dropped its trailing Cohort
# clause - the query now ends after the "Response to Treatment" line
# with no trailing comma / newline.
$newCasesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
"WHERE demo.breed IN ['Parson Russell Terrier'] `n" +
"MATCH (c)<--(diag:diagnosis)`n" +
"OPTIONAL MATCH (samp:sample)-->(c)`n" +
"OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" +
"WITH DISTINCT c, s, demo, diag, co`n" +
"RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" +
"        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" +
"        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" +
"        coalesce(demo.breed, '') AS Breed ,`n" +
"        coalesce(diag.disease_term, '') AS Diagnosis ,`n" +
"        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" +
"        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" +
"        coalesce(demo.sex, '') AS Sex ,`n" +
"        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" +
"        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" +
"        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $newCasesQuery

# Selection ends up on B2 after the edit, with the view scrolled back
# to the top of the sheet.
$ws.Range("B2").Select()
